$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.404.24"
$ws.Range("E2").Value = "  -1.67%  "

$ws.Range("D3").Value = "3.389.49"
$ws.Range("E3").Value = "  -1.78%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "404.19"
$ws.Range("E5").Value = "  -0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.20"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.588"
$ws.Range("E7").Value = "  -3.95%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.680"
$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.127"
$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.45"
$ws.Range("E11").Value = "  -3.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.140"
$ws.Range("E12").Value = "  -1.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.68"
$ws.Range("E13").Value = "  -1.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.31"
$ws.Range("E14").Value = "  -4.35%  "

$ws.Range("D15").Value = "3.412.09"
$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.63"
$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").Value = "61.446.21"
$ws.Range("E17").Value = "  -1.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.00"
$ws.Range("E18").Value = "  -3.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000141"
$ws.Range("E19").Value = "  +2.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("E20").Value = "  -4.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "82.90"
$ws.Range("E21").Value = "  -0.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "309.91"
$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.66"
$ws.Range("E23").Value = "  -3.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.14"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.79"
$ws.Range("E25").Value = "  +9.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.18"
$ws.Range("E26").Value = "  +8.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.34"
$ws.Range("E27").Value = "  -2.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.01"
$ws.Range("E28").Value = "  -5.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.70"
$ws.Range("E29").Value = "  +4.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "43.75"
$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("E31").Value = "  -2.15%  "

$ws.Range("E32").Value = "  -2.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.25"
$ws.Range("E33").Value = "  -4.07%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0479"
$ws.Range("E35").Value = "  -1.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.26"
$ws.Range("E36").Value = "  -2.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").Value = "  -4.56%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.317"
$ws.Range("E40").Value = "  +11.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "139.60"
$ws.Range("E41").Value = "  +1.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.124"
$ws.Range("E42").Value = "  -2.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.95"
$ws.Range("E43").Value = "  -1.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.91"
$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.60"
$ws.Range("E45").Value = "  -3.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("E46").Value = "  -1.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.94"
$ws.Range("E47").Value = "  -3.99%  "

$ws.Range("D48").Value = "2.087.75"
$ws.Range("E48").Value = "  -3.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.32"
$ws.Range("E49").Value = "  -2.23%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.89"
$ws.Range("E50").Value = "  +1.74%  "

$ws.Range("B51").Value = "Fetch.AI"
$ws.Range("C51").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("E51").Value = "  +15.27%  "
